# User can add items by entering data into the console
# Adds two new columns (Sold Price, Is sold), renames existing headers,
# and appends the new field values for the existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename / repurpose existing headers ---
$ws.Range("C1").Value = "Bought price"
$ws.Range("D1").Value = "Sold Price"

# --- New headers ---
$ws.Range("E1").Value = "Condition"
$ws.Range("F1").Value = "Is sold"

# --- Update the existing data row ---
# Keep the date as literal text (not auto-converted to a date serial number)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "01-02-2025"
$ws.Range("A2").ClearFormats()

$ws.Range("C2").Value = 13.39
$ws.Range("D2").Value = 14.56
$ws.Range("E2").Value = "New With Tag"
$ws.Range("F2").Value = $true

# --- Resize columns to fit the new content ---
$ws.Columns("A:F").AutoFit()
